# Apply updated crawl data (想去人数 / 地点 / Cover) to the "展览" and
# "全部类型" worksheets of the Suzhou comic-convention workbook.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 141
    $ws.Range("F3").Value = 1668
    $ws.Range("F5").Value = 1112
    $ws.Range("F7").Value = 11743

    $ws.Range("D11").Value = "兴中路与鲈乡北路交汇处 香漫商业广场"
    $ws.Range("F11").Value = 388
    $ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202403/DI2ackIO1710137864319.jpeg"

    $ws.Range("F13").Value = 830
    $ws.Range("F14").Value = 13417
    $ws.Range("F15").Value = 13278
    $ws.Range("F23").Value = 149
}

$wb.Save()
